$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 11:33"

# Row 16 - Belgica
$ws.Range("B16").Value = 50267
$ws.Range("C16").Value = 361
$ws.Range("D16").Value = 12378
$ws.Range("E16").Value = 29965
$ws.Range("F16").Value = 655
$ws.Range("G16").Value = 80
$ws.Range("H16").Value = 7924

# Row 30 - Bielorrusia
$ws.Range("B30").Value = 17489
$ws.Range("C30").Value = 784
$ws.Range("D30").Value = 3259
$ws.Range("E30").Value = 14127
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 103

# Row 39 - Indonesia
$ws.Range("B39").Value = 11587
$ws.Range("C39").Value = 395
$ws.Range("D39").Value = 1954
$ws.Range("E39").Value = 8769
$ws.Range("G39").Value = 19
$ws.Range("H39").Value = 864

# Row 54 - Finlandia
$ws.Range("B54").Value = 5327
$ws.Range("C54").Value = 73
$ws.Range("E54").Value = 2097

# Row 128 - Montenegro
$ws.Range("B128").Value = 323
$ws.Range("C128").Value = 1
$ws.Range("E128").Value = 66
